$d = $word.ActiveDocument

# Locate the "Platform impact..." bullet (last bullet under KEY ACHIEVEMENTS
# AND IMPACT, immediately before the "TECHNICAL SKILLS" heading) and insert
# four new achievement bullets after it, restoring content that had been
# truncated by pipe-delimited bullet expansion (total_max 5 -> 9).

$anchorSnippet = "Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"

$anchorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains($anchorSnippet)) {
        $anchorPara = $p
    }
}

$newBullets = @(
    "• Real-time collaboration at national scale",
    "• Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from ±4.2% to ±2.1%",
    "• Increased voter turnout prediction accuracy from 71% to 87%",
    "• Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis"
)

$r = $anchorPara.Range
$insertIndex = $anchorPara.Index

foreach ($bulletText in $newBullets) {
    $r.InsertParagraphAfter()
    $insertIndex = $insertIndex + 1
    $newPara = $d.Paragraphs.Item($insertIndex)
    $newPara.Range.Text = $bulletText
    $r = $newPara.Range
}
